# Add normalized cell counts (per Grant's request) to the first 23 rows
# of columns B:G, replacing the raw counts with values normalized against
# an external reference (Cell Count Data (Final MATLAB Format) v2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalized replacement values for rows 1-23, columns B,C,D,E,F,G
$newValues = @(
    @(1609.105180327869, 1711.8143015521064, 1483.1827309236949, 1852.7306873669716, 1737.7751312071055, 1735.2635869565217),
    @(1720.5112131147541, 1804.7311529933479, 1494.2006597819852, 1913.7222987354453, 1763.3181267662496, 1691.2345108695654),
    @(1764.906098360656, 1835.1078159645231, 1564.5458978772233, 1950.8476273945164, 1767.5752926927737, 1753.3932065217393),
    @(1861.2346229508198, 1827.9603658536585, 1616.2454102122776, 1995.0444472267438, 1861.232943076302, 1838.8614130434785),
    @(1908.1424262295084, 1910.1560421286028, 1679.8103843947217, 2094.9292600475778, 1886.7759386354462, 1884.6171195652175),
    @(1988.5558032786887, 2005.7531873614189, 1775.5816121629375, 2192.1622636784778, 1980.4335890189748, 1970.948641304348),
    @(2096.6112786885246, 2086.1620011086475, 1824.7385255306942, 2288.5113309127337, 2036.6281792490918, 2082.3163043478262),
    @(2200.4785573770496, 2170.1445399113081, 1903.559093516925, 2359.2262426442971, 2146.4630601534113, 2123.755434782609),
    @(2260.7885901639347, 2220.1766906873613, 2003.567986230637, 2505.0757480906477, 2226.4977795720633, 2203.1804347826087),
    @(2356.2794754098363, 2377.4205931263855, 2114.594807802639, 2594.3533241517471, 2323.5611626968107, 2330.951086956522),
    @(2618.4605901639347, 2608.8192904656316, 2303.594664371773, 2863.0699887316896, 2554.2995559144128, 2532.1035326086958),
    @(2719.8149508196725, 2724.0719235033257, 2461.2358003442341, 3036.3215224740206, 2664.1344368187324, 2645.1978260869569),
    @(2907.4461639344263, 2880.4223946784919, 2590.060814687321, 3197.1979466633284, 2837.8268066209125, 2832.5372282608696),
    @(3098.42793442623, 3076.9772727272725, 2778.2131382673551, 3406.6908726680863, 3012.3706096083974, 3060.452445652174),
    @(3239.1513442622954, 3304.8022450110861, 2913.818416523236, 3509.2274946788539, 3193.7258780783209, 3209.8059782608698),
    @(3237.4760655737709, 3150.238636363636, 2935.0067412507174, 3517.1829222486544, 3228.6346386758178, 3166.6402173913048),
    @(3278.5203934426231, 3322.6708702882479, 3024.8452380952381, 3634.7464630023796, 3263.5433992733147, 3252.1084239130437),
    @(3260.0923278688529, 3285.1467572062083, 2997.7241824440621, 3595.8532615500194, 3226.0803391199033, 3206.3527173913044),
    @(3320.402360655738, 3274.4255820399112, 3018.9125071715434, 3646.2376361587585, 3245.663302381914, 3235.7054347826088),
    @(3272.656918032787, 3215.4591186252769, 2985.8587205966724, 3599.3890071365977, 3198.8344771901498, 3220.1657608695655),
    @(3254.2288524590167, 3199.3773558758312, 2948.5672690763054, 3543.7010141479909, 3201.3887767460642, 3158.8703804347829),
    @(3219.8856393442625, 3255.663525498891, 2931.6166092943204, 3515.4150494553655, 3164.7771497779572, 3135.5608695652177),
    @(3202.2952131147545, 3234.2211751662967, 2946.8722030981066, 3525.1383498184555, 3163.9257165926529, 3163.1869565217394)
)

# Column B is spreadsheet column index 2, through G = column index 7
$colIndexes = @(2, 3, 4, 5, 6, 7)

for ($r = 0; $r -lt $newValues.Length; $r++) {
    $rowValues = $newValues[$r]
    for ($c = 0; $c -lt $colIndexes.Length; $c++) {
        $ws.Cells.Item($r + 1, $colIndexes[$c]).Value = $rowValues[$c]
    }
}

# Update the selection to match the newly-edited range (B1:G23), as the
# author left the selection here after entering the normalized data.
$ws.Range("B1:G23").Select()
